$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recalculated "K" column (strike count, replacing legacy Strike# values)
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 3
$ws.Range("G8").Value = 2
$ws.Range("G9").Value = 1
